$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.424.67'
$ws.Range("E2").Value = '  -4.21%  '
$ws.Range("D3").Value = '2.394.91'
$ws.Range("E3").Value = '  -4.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '499.70'
$ws.Range("E5").Value = '  -6.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.83'
$ws.Range("E6").Value = '  -3.72%  '
$ws.Range("E7").Value = '  -0.39%  '
$ws.Range("E8").Value = '  -3.28%  '
$ws.Range("D9").Value = '2.391.46'
$ws.Range("E9").Value = '  -4.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0956'
$ws.Range("E10").Value = '  -3.66%  '
$ws.Range("E11").Value = '  -1.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.319'
$ws.Range("E12").Value = '  -3.15%  '
$ws.Range("E13").Value = '  -9.79%  '
$ws.Range("D14").Value = '2.817.91'
$ws.Range("E14").Value = '  -4.38%  '
$ws.Range("D15").Value = '57.009.07'
$ws.Range("E15").Value = '  -2.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.47'
$ws.Range("E16").Value = '  -3.74%  '
$ws.Range("E17").Value = '  -3.51%  '
$ws.Range("D18").Value = '2.341.74'
$ws.Range("E18").Value = '  -6.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.14'
$ws.Range("E19").Value = '  -4.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '311.81'
$ws.Range("E21").Value = '  -5.18%  '
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.53'
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("D26").Value = '2.488.62'
$ws.Range("E26").Value = '  -5.11%  '
$ws.Range("E27").Value = '  -9.40%  '
$ws.Range("E28").Value = '  -6.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.19'
$ws.Range("E29").Value = '  -3.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.07'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.65'
$ws.Range("E31").Value = '  -4.59%  '
$ws.Range("D32").Value = '0.0₃0710'
$ws.Range("E32").Value = '  -6.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.09'
$ws.Range("E33").Value = '  -2.72%  '
$ws.Range("E34").Value = '  -6.55%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.73'
$ws.Range("E37").Value = '  -1.88%  '
$ws.Range("E38").Value = '  -0.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.73'
$ws.Range("E39").Value = '  -5.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35.88'
$ws.Range("E40").Value = '  -1.39%  '
$ws.Range("E41").Value = '  -6.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.771'
$ws.Range("E42").Value = '  -6.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '129.36'
$ws.Range("E43").Value = '  -1.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.34'
$ws.Range("E44").Value = '  -3.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.85'
$ws.Range("E45").Value = '  -3.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.570'
$ws.Range("E46").Value = '  -3.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '252.97'
$ws.Range("E47").Value = '  -7.92%  '
$ws.Range("E48").Value = '  -3.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0484'
$ws.Range("E49").Value = '  -4.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.81'
$ws.Range("E50").Value = '  -4.34%  '
$ws.Range("E51").Value = '  -5.07%  '
